# Updates cryptos list price (D) and 1h volume change (E) columns
# to match the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.688.31'
$ws.Range("E2").Value = '  +3.44%  '

$ws.Range("D3").Value = '1.863.44'
$ws.Range("E3").Value = '  +2.87%  '

$ws.Range("E4").Value = '  +0.42%  '

$ws.Range("D5").Value = '''231.22'
$ws.Range("E5").Value = '  +2.48%  '

$ws.Range("D8").Value = '''42.51'
$ws.Range("E8").Value = '  +11.13%  '

$ws.Range("E9").Value = '  +7.21%  '

$ws.Range("D10").Value = '''0.0696'
$ws.Range("E10").Value = '  +3.02%  '

$ws.Range("E11").Value = '  +3.95%  '

$ws.Range("D12").Value = '2.133.89'
$ws.Range("E12").Value = '  +2.94%  '

$ws.Range("D13").Value = '''11.64'
$ws.Range("E13").Value = '  +3.56%  '

$ws.Range("D14").Value = '1.863.84'
$ws.Range("E14").Value = '  +2.98%  '

$ws.Range("D15").Value = '''0.680'
$ws.Range("E15").Value = '  +7.32%  '

$ws.Range("D16").Value = '''4.73'
$ws.Range("E16").Value = '  +6.83%  '

$ws.Range("D17").Value = '35.707.73'
$ws.Range("E17").Value = '  +3.59%  '

$ws.Range("D18").Value = '''70.46'
$ws.Range("E18").Value = '  +3.07%  '

$ws.Range("D19").Value = '''248.98'
$ws.Range("E19").Value = '  +2.32%  '

$ws.Range("D20").Value = '0.0₃0805'
$ws.Range("E20").Value = '  +4.03%  '

$ws.Range("D21").Value = '''12.30'
$ws.Range("E21").Value = '  +9.59%  '

$ws.Range("D22").Value = '''4.75'
$ws.Range("E22").Value = '  +15.01%  '

$ws.Range("E23").Value = '  +0.33%  '

$ws.Range("D24").Value = '''2.23'
$ws.Range("E24").Value = '  +0.90%  '

$ws.Range("D25").Value = '''170.63'
$ws.Range("E25").Value = '  -0.07%  '

$ws.Range("D26").Value = '''8.00'
$ws.Range("E26").Value = '  +2.87%  '

$ws.Range("D27").Value = '''17.92'
$ws.Range("E27").Value = '  +1.37%  '

$ws.Range("E28").Value = '  +1.59%  '

$ws.Range("E29").Value = '  +16.78%  '

$ws.Range("E30").Value = '  +0.49%  '

$ws.Range("D31").Value = '3.319.07'
$ws.Range("E31").Value = '  +36.60%  '

$ws.Range("E32").Value = '  +6.06%  '

$ws.Range("D33").Value = '''0.0547'
$ws.Range("E33").Value = '  +5.86%  '

$ws.Range("D34").Value = '''3.96'
$ws.Range("E34").Value = '  +4.18%  '

$ws.Range("D35").Value = '''1.90'
$ws.Range("E35").Value = '  +4.23%  '

$ws.Range("D36").Value = '''101.67'
$ws.Range("E36").Value = '  +24.08%  '

$ws.Range("D37").Value = '''0.691'
$ws.Range("E37").Value = '  +7.32%  '

$ws.Range("D38").Value = '1.370.76'
$ws.Range("E38").Value = '  +1.01%  '

$ws.Range("D39").Value = '''2.49'
$ws.Range("E39").Value = '  +6.41%  '

$ws.Range("D40").Value = '''1.09'
$ws.Range("E40").Value = '  +3.43%  '

$ws.Range("E41").Value = '  +4.91%  '

$ws.Range("D42").Value = '''1.01'
$ws.Range("E42").Value = '  +5.81%  '

$ws.Range("D43").Value = '''1.26'
$ws.Range("E43").Value = '  +4.13%  '

$ws.Range("D44").Value = '''14.87'
$ws.Range("E44").Value = '  +8.04%  '

$ws.Range("E45").Value = '  +1.10%  '

$ws.Range("E46").Value = '  +1.26%  '

$ws.Range("D47").Value = '''6.29'
$ws.Range("E47").Value = '  +8.67%  '

$ws.Range("D48").Value = '''0.0522'
$ws.Range("E48").Value = '  +2.50%  '

$ws.Range("D49").Value = '2.032.41'
$ws.Range("E49").Value = '  +2.96%  '

$ws.Range("D50").Value = '''104.69'
$ws.Range("E50").Value = '  +1.98%  '

$ws.Range("E51").Value = '  +0.39%  '

